$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B20 to be a real number instead of text
$ws.Range("B20").Value = 21

# Add row 21
$ws.Range("A21").Value = "Vanda Dyy"
$ws.Range("B21").Value = 23
$ws.Range("C21").Value = "Male"
$ws.Range("D21").Value = "Phnom Penh"
$ws.Range("E21").Value = "Class b 2025"
$ws.Range("F21").Value = "image\44c795a0026549cea99c8f4d0d600342.png"

# Add row 22
$ws.Range("A22").Value = "dyy"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "19"
$ws.Range("C22").Value = "Male"
$ws.Range("D22").Value = "KPC"
$ws.Range("E22").Value = "B2025"
$ws.Range("F22").Value = "image\0b309e35ab6d40738af04a70c6525f40.png"
